$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the two new worksheets (view_centers, edit_centers) at the end of
#    the workbook, after add_new_centers.
# ---------------------------------------------------------------------------
$afterSheet1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$viewCenters = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet1)
$viewCenters.Name = "view_centers"

$afterSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$editCenters = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet2)
$editCenters.Name = "edit_centers"

# ---------------------------------------------------------------------------
# 2) Populate view_centers (row / runmode table used for the "search center"
#    test cases).
# ---------------------------------------------------------------------------
$viewCenters.Range("A1:A5").NumberFormat = "@"

$viewCenters.Range("A1").Value = "row"
$viewCenters.Range("B1").Value = "runmode"
$viewCenters.Range("A2").Value = "3"
$viewCenters.Range("B2").Value = "Y"
$viewCenters.Range("A3").Value = "4"
$viewCenters.Range("B3").Value = "Y"
$viewCenters.Range("A4").Value = "5"
$viewCenters.Range("B4").Value = "Y"
$viewCenters.Range("A5").Value = "6"
$viewCenters.Range("B5").Value = "Y"

$viewCenters.Columns.Item(1).ColumnWidth = 8.307291666666666
$viewCenters.PageSetup.Orientation = 1
$viewCenters.Range("A1:B5").Select()

# ---------------------------------------------------------------------------
# 3) Populate edit_centers (row / new code / new name / runmode table used
#    for the "edit center" test cases).
# ---------------------------------------------------------------------------
$editCenters.Range("A1:C5").NumberFormat = "@"

$editCenters.Range("A1").Value = "row"
$editCenters.Range("B1").Value = "new code"
$editCenters.Range("C1").Value = "new name"
$editCenters.Range("D1").Value = "runmode"

$editCenters.Range("A2").Value = "3"
$editCenters.Range("B2").Value = "MCEN1"
$editCenters.Range("D2").Value = "Y"
$editCenters.Range("A3").Value = "4"
$editCenters.Range("B3").Value = "MCEN2"
$editCenters.Range("D3").Value = "Y"
$editCenters.Range("A4").Value = "5"
$editCenters.Range("B4").Value = "MCEN3"
$editCenters.Range("D4").Value = "Y"
$editCenters.Range("A5").Value = "6"
$editCenters.Range("B5").Value = "MCEN4"
$editCenters.Range("D5").Value = "Y"

$editCenters.Range("C2").Value = "NewName1"
$editCenters.Range("C3").Value = "NewName2"
$editCenters.Range("C4").Value = "NewName3"
$editCenters.Range("C5").Value = "NewName4"

$editCenters.Columns.Item(3).ColumnWidth = 11.833333333333332
$editCenters.Range("B2").Select()

Write-Host "done"
